$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tresquintos.cl")

# Insert a new row at position 20 (shifts rows 20:48 down to 21:49,
# carrying their existing formatting with them, and the new row 20
# inherits the default row-2 style already used throughout the table).
$ws.Rows.Item(20).Insert()

# Populate the freshly inserted row with the new candidate record.
$ws.Cells.Item(20, 1).Value = 48
$ws.Cells.Item(20, 2).Value = 15
$ws.Cells.Item(20, 3).Value = 5.5
$ws.Cells.Item(20, 4).Value = "Metropolitana"
$ws.Cells.Item(20, 5).Value = "Nathalie Joignant"
$ws.Cells.Item(20, 6).Value = "Partido Ecologista Verde"
$ws.Cells.Item(20, 7).Value = "PEV"
$ws.Cells.Item(20, 8).Value = 1

# Move the active selection the way it ended up after the edit.
$ws.Range("E12").Select()

# The sheet's hidden _FilterDatabase name (backing the AutoFilter) and the
# worksheet-level remembered sort range both need to grow by one row too.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $item = $names.Item($i)
    if ($item.Name -eq "tresquintos.cl!_FilterDatabase") {
        $item.RefersTo = "=tresquintos.cl!`$A`$1:`$H`$45"
    }
}

$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("C2:C50"))
$sort.SortFields.Add($ws.Range("F2:F50"))
$sort.SetRange($ws.Range("A2:H50"))
$sort.Header = 2
$sort.Apply()
